$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 333
$ws.Range("F3").Value = 537
$ws.Range("F4").Value = 675
$ws.Range("F6").Value = 890
$ws.Range("F7").Value = 424
$ws.Range("F8").Value = 115
$ws.Range("F9").Value = 432
$ws.Range("F12").Value = 1137
$ws.Range("F14").Value = 45
$ws.Range("F15").Value = 1956
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 29
$ws.Range("F20").Value = 514
$ws.Range("F23").Value = 546
$ws.Range("F24").Value = 375
$ws.Range("F25").Value = 375
$ws.Range("F26").Value = 692
$ws.Range("F27").Value = 443
$ws.Range("F28").Value = 2725
$ws.Range("F31").Value = 3195
$ws.Range("F32").Value = 637
$ws.Range("F33").Value = 513
$ws.Range("F34").Value = 222
$ws.Range("F35").Value = 965
$ws.Range("F36").Value = 725
$ws.Range("F38").Value = 638
$ws.Range("F39").Value = 620

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 62
$ws.Range("F22").Value = 170
$ws.Range("F23").Value = 128
$ws.Range("F24").Value = 443

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 2932
$ws.Range("F5").Value = 253
$ws.Range("F6").Value = 369

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 537
$ws.Range("F6").Value = 253
$ws.Range("F7").Value = 675
$ws.Range("F9").Value = 890
$ws.Range("F10").Value = 424
$ws.Range("F11").Value = 115
$ws.Range("F12").Value = 432
$ws.Range("F16").Value = 1137
$ws.Range("F19").Value = 45
$ws.Range("F20").Value = 369
$ws.Range("F21").Value = 1957
$ws.Range("F22").Value = 1958
$ws.Range("F25").Value = 62
$ws.Range("F27").Value = 514
$ws.Range("F34").Value = 375
$ws.Range("F36").Value = 692
$ws.Range("F37").Value = 443
$ws.Range("F39").Value = 2725
$ws.Range("F41").Value = 3195
$ws.Range("F42").Value = 637
$ws.Range("F43").Value = 222
$ws.Range("F44").Value = 965
$ws.Range("F47").Value = 128
$ws.Range("F48").Value = 443
$ws.Range("F49").Value = 725
$ws.Range("F51").Value = 638
$ws.Range("F52").Value = 620
